# Auto-generated edit script: updates crypto price/volume table
# to reflect the "Sat May 13 18:40:19 UTC 2023" GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.503.20'

$ws.Range("D3").Value = '1.841.28'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.028'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = '  +2.64%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.80'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +4.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.024'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  +2.34%  '

$ws.Range("E7").Value = '  +3.26%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3728'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = '  +3.44%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07391'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  +3.47%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8762'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  +4.58%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.48'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  +5.16%  '

$ws.Range("D12").Value = '1.851.70'

$ws.Range("E12").Value = '  +5.23%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.496'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  +4.60%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.677'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +3.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07140'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  +3.50%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.65'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +4.61%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.031'
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009022'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  +4.24%  '

$ws.Range("E19").Value = '  +2.48%  '

$ws.Range("E20").Value = '  +3.35%  '

$ws.Range("D21").Value = '27.518.17'

$ws.Range("E21").Value = '  +4.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.232'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  +2.69%  '

$ws.Range("E23").Value = '  +2.60%  '

$ws.Range("D24").Value = '2.067.04'

$ws.Range("E24").Value = '  +4.44%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.93'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  +3.53%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.921'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  +8.14%  '

$ws.Range("E27").Value = '  +3.94%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.256'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  +3.66%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.940'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  +5.83%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.29'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  +1.87%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09068'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  +2.64%  '

$ws.Range("E32").Value = '  +7.70%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7652'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +5.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.488'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  +4.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.873'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  +5.19%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.027'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  +2.75%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.145'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  +5.17%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01970'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  +4.48%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05252'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +2.80%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5175'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  +5.11%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.779'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  +6.92%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1664'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  +3.36%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.623'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  +4.74%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.526'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  +5.58%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '109.06'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  +4.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.55'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  +3.70%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.028'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  +2.73%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.705'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  +4.87%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4643'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +4.48%  '

$ws.Range("B50").Value = 'RenderToken'

$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.899'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  +10.84%  '

$ws.Range("B51").Value = 'Cronos'

$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06329'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  +2.54%  '
